# Update "想去人数" (number of people interested) counts on the
# "展览" and "全部类型" sheets to reflect the latest scrape.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 105
    $ws.Range("F5").Value = 2808
    $ws.Range("F6").Value = 276
}
